# Cost_Savings_Summary.xlsx -- finished Word document report
#
# The underlying simulation was rerun with a corrected unit-cost scale
# (cost figures come back ~1e4-1e5x smaller than before) and a slightly
# different scenario_qty/scenario_value draw. savings_delta/percent and
# the Totals sheet are refreshed to match, and the now-redundant
# "note" column on the Totals sheet (explaining which subprocesses were
# excluded) is dropped.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1 (Savings_by_subprocess): recalculated cost/savings figures ---
$ws1.Range("C2").Value = 26694.981
$ws1.Range("D2").Value = 48119.94
$ws1.Range("E2").Value = 24541.1694
$ws1.Range("F2").Value = 2153.811599999997
$ws1.Range("G2").Value = 0.08068226757681518
$ws1.Range("C3").Value = 4399.7565
$ws1.Range("D3").Value = 114055.85
$ws1.Range("E3").Value = 4106.0106
$ws1.Range("F3").Value = 293.7458999999999
$ws1.Range("G3").Value = 0.06676412660564282
$ws1.Range("C4").Value = 181.4614
$ws1.Range("D4").Value = 52475.5
$ws1.Range("E4").Value = 199.4069
$ws1.Range("F4").Value = -17.94550000000001
$ws1.Range("G4").Value = -0.09889431030511178
$ws1.Range("C5").Value = 26694.981
$ws1.Range("D5").Value = 47754.77166666667
$ws1.Range("E5").Value = 24354.93355
$ws1.Range("F5").Value = 2340.047449999998
$ws1.Range("G5").Value = 0.08765870445833987
$ws1.Range("C6").Value = 4399.7565
$ws1.Range("D6").Value = 113348.92
$ws1.Range("E6").Value = 4080.56112
$ws1.Range("F6").Value = 319.1953799999997
$ws1.Range("G6").Value = 0.07254841944093944
$ws1.Range("C7").Value = 181.4614
$ws1.Range("D7").Value = 50567.3
$ws1.Range("E7").Value = 192.15574
$ws1.Range("F7").Value = -10.69433999999998
$ws1.Range("G7").Value = -0.05893451720310756
$ws1.Range("C8").Value = 26694.981
$ws1.Range("D8").Value = 47545.58583333333
$ws1.Range("E8").Value = 24248.248775
$ws1.Range("F8").Value = 2446.732225
$ws1.Range("G8").Value = 0.09165514015537227
$ws1.Range("C9").Value = 4399.7565
$ws1.Range("D9").Value = 112944.96
$ws1.Range("E9").Value = 4066.01856
$ws1.Range("F9").Value = 333.73794
$ws1.Range("G9").Value = 0.07585372963253763
$ws1.Range("C10").Value = 181.4614
$ws1.Range("D10").Value = 49136.15000000001
$ws1.Range("E10").Value = 186.71737
$ws1.Range("F10").Value = -5.255970000000048
$ws1.Range("G10").Value = -0.02896467237660487

# --- Sheet2 (KPI_changes): recalculated scenario_value / change figures ---
$ws2.Range("C2").Value = 20990.58638888889
$ws2.Range("D2").Value = 20711.85230442368
$ws2.Range("C3").Value = 48119.94
$ws2.Range("D3").Value = -4223.159999999996
$ws2.Range("C4").Value = 114055.85
$ws2.Range("D4").Value = -8159.608333333323
$ws2.Range("C5").Value = 52475.5
$ws2.Range("D5").Value = 4722.5
$ws2.Range("C6").Value = 23614.4
$ws2.Range("D6").Value = 23335.66591553479
$ws2.Range("C7").Value = 47754.77166666667
$ws2.Range("D7").Value = -4588.328333333331
$ws2.Range("C8").Value = 113348.92
$ws2.Range("D8").Value = -8866.53833333333
$ws2.Range("C9").Value = 50567.3
$ws2.Range("D9").Value = 2814.299999999996
$ws2.Range("C10").Value = 24926.31166666667
$ws2.Range("D10").Value = 24647.57758220146
$ws2.Range("C11").Value = 47545.58583333333
$ws2.Range("D11").Value = -4797.514166666668
$ws2.Range("C12").Value = 112944.96
$ws2.Range("D12").Value = -9270.498333333337
$ws2.Range("C13").Value = 49136.15000000001
$ws2.Range("D13").Value = 1383.150000000009

# --- Sheet3 (Totals): drop the "note" column, recalc totals ---
$ws3.Columns.Item(5).Delete()
$ws3.Range("A2").Value = 31276.1989
$ws3.Range("B2").Value = 28846.5869
$ws3.Range("C2").Value = 2429.611999999997
$ws3.Range("D2").Value = 0.07768245776183491
$ws3.Range("A3").Value = 31276.1989
$ws3.Range("B3").Value = 28627.65041
$ws3.Range("C3").Value = 2648.548489999998
$ws3.Range("D3").Value = 0.08468255680520045
$ws3.Range("A4").Value = 31276.1989
$ws3.Range("B4").Value = 28500.984705
$ws3.Range("C4").Value = 2775.214195
$ws3.Range("D4").Value = 0.08873246406551021

# --- Active sheet/tab switches from Totals back to Savings_by_subprocess ---
$ws1.Activate()
